$wb = $excel.ActiveWorkbook

# Sheet "展览" and "全部类型" both have the same F-column updates applied.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 320
    $ws.Range("F3").Value = 75
    $ws.Range("F5").Value = 4696
    $ws.Range("F6").Value = 368
    $ws.Range("F9").Value = 729

    if ($name -eq "展览") {
        $ws.Range("F10").Value = 208
    } else {
        $ws.Range("F11").Value = 208
    }
}
